$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column D entirely (S2 scenario), shifting nothing else since it's the last column
$ws.Range("D1:D6").Delete()

# Update B2:C2 values
$ws.Range("B2").Value = 0.984
$ws.Range("C2").Value = 0.016

# Update B3:C6 values (rows 3 through 6 all become 0.96 / 0.04)
$ws.Range("B3:B6").Value = 0.96
$ws.Range("C3:C6").Value = 0.04

# Add new rows 7 and 8
$ws.Range("A7").Value = 2045
$ws.Range("B7").Value = 0.96
$ws.Range("C7").Value = 0.04

$ws.Range("A8").Value = 2050
$ws.Range("B8").Value = 0.96
$ws.Range("C8").Value = 0.04
